{"js": "// Resume update: refresh \"TECHNICAL SKILLS\" tech stacks, retitle the\n// flagship product, and tweak two bullet/interest lines.\n//\n// Each change is applied as a single search-and-replace against the\n// document body. Word's search matches across run boundaries, and\n// insertText(..., Word.InsertLocation.replace) collapses the match\n// into one run (inheriting the formatting of the first run of the\n// match), which is exactly what the target edit does for these runs\n// (they all share the same sz/szCs 22 formatting).\n\nconst replacements = [\n  // 1) Development stack: drop Webpack/Karma/Shell, add Express/Knex.js/\n  //    TypeScript, reorder Bootstrap/D3.js, and shorten CSS3 -> CSS.\n  [\n    \"Angular, Node.js, Webpack, Karma, Bootstrap, D3.js, JavaScript, jQuery, Sass, LESS, CSS3, HTML5, Shell\",\n    \"Angular, Node.js, Express, Knex.js, D3.js, Bootstrap, TypeScript, JavaScript, jQuery, Sass, LESS, CSS, HTML5\",\n  ],\n  // 2) Systems stack: drop Solr/Redis/Nginx/Jenkins, add AWS/Docker/PostgreSQL.\n  [\n    \"MySQL, SQLite, Solr, Redis, Nginx, Jenkins\",\n    \"AWS, Docker, PostgreSQL, MySQL, SQLite\",\n  ],\n  // 3) Product rename SmartView -> ActiveInsight, and a small wording tweak.\n  [\n    \"Designed and wrote functionality for SmartView, an emotional analytics SaaS platform which provides enterprise solutions through brand equity assessment, market landscaping, competitive intelligence, and deep insights into different markets\",\n    \"Designed and wrote functionality for ActiveInsight, an emotional analytics SaaS platform which provides enterprise solutions through brand equity assessment, market landscaping, competitive intelligence, and deep market insights\",\n  ],\n  // 4) \"leveraging\" -> \"with\", and append \", and D3.js\" to the tooling list.\n  [\n    \"Architected and implemented the front-end web client leveraging modern tooling- Node.js, Karma, Webpack, and Angular\",\n    \"Architected and implemented the front-end web client with modern tooling- Node.js, Karma, Webpack, Angular, and D3.js\",\n  ],\n  // 5) Interests: swap \"Graphic Design\" for \"Piano, Design\".\n  [\n    \", Graphic Design\",\n    \", Piano, Design\",\n  ],\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Resume update: refresh \"TECHNICAL SKILLS\" tech stacks, retitle the\n# flagship product, and tweak two bullet/interest lines.\n#\n# Each change is a single Find/Replace (wdReplaceAll, restricted to the\n# one match that exists) over $d.Content. Word's Find matches across run\n# boundaries and the matched range collapses into a single run carrying\n# the formatting of the first run in the match, matching the target edit\n# (all affected runs share the same sz/szCs 22 formatting).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) Development stack: drop Webpack/Karma/Shell, add Express/Knex.js/\n#    TypeScript, reorder Bootstrap/D3.js, and shorten CSS3 -> CSS.\nReplace-Text \"Angular, Node.js, Webpack, Karma, Bootstrap, D3.js, JavaScript, jQuery, Sass, LESS, CSS3, HTML5, Shell\" \"Angular, Node.js, Express, Knex.js, D3.js, Bootstrap, TypeScript, JavaScript, jQuery, Sass, LESS, CSS, HTML5\"\n\n# 2) Systems stack: drop Solr/Redis/Nginx/Jenkins, add AWS/Docker/PostgreSQL.\nReplace-Text \"MySQL, SQLite, Solr, Redis, Nginx, Jenkins\" \"AWS, Docker, PostgreSQL, MySQL, SQLite\"\n\n# 3) Product rename SmartView -> ActiveInsight, and a small wording tweak.\nReplace-Text \"Designed and wrote functionality for SmartView, an emotional analytics SaaS platform which provides enterprise solutions through brand equity assessment, market landscaping, competitive intelligence, and deep insights into different markets\" \"Designed and wrote functionality for ActiveInsight, an emotional analytics SaaS platform which provides enterprise solutions through brand equity assessment, market landscaping, competitive intelligence, and deep market insights\"\n\n# 4) \"leveraging\" -> \"with\", and append \", and D3.js\" to the tooling list.\nReplace-Text \"Architected and implemented the front-end web client leveraging modern tooling- Node.js, Karma, Webpack, and Angular\" \"Architected and implemented the front-end web client with modern tooling- Node.js, Karma, Webpack, Angular, and D3.js\"\n\n# 5) Interests: swap \"Graphic Design\" for \"Piano, Design\".\nReplace-Text \", Graphic Design\" \", Piano, Design\"\n"}
